# Savesheet para diferentes bimestres
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "Ernane-Des. Tec. M"

# Row 3
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "Ernane-Des. Tec. M"

# Row 4
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "João Bosco-Gestão Integr"
$ws.Range("F4").Value = "Ernane-Des. Tec. M"

# Row 6
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "João Bosco-Gestão Integr"
$ws.Range("F6").Value = "Ernane-Des. Tec. M"

# Row 7
$ws.Range("B7").Value = "Aline S. M.-T. M. Metali"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "Ernane-Des. Tec. M"

# Row 8
$ws.Range("B8").Value = "Aline S. M.-T. M. Metali"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "Ernane-Des. Tec. M"
